$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new daily observations appended below the existing data
# (05-10-2021 was previously the last row -> 192).
$startRow = 193

$dates = @("06-10-2021", "07-10-2021", "08-10-2021")
$tcm   = @(117.53, 118.07, 117.82)
$tcm5  = @(192.06, 193.05, 192.54)
$tcmx  = @(103.62, 104.04, 103.85)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i

    # Write the date as a text formula (quoted literal) so Excel stores it
    # as a plain string instead of auto-converting "dd-mm-yyyy" looking
    # text into a date serial number / applying a date number format.
    $ws.Cells.Item($r, 1).Formula = "=""" + $dates[$i] + """"

    $ws.Cells.Item($r, 2).Value = $tcm[$i]
    $ws.Cells.Item($r, 3).Value = $tcm5[$i]
    $ws.Cells.Item($r, 4).Value = $tcmx[$i]
}

# Convert the formulas in column A to their plain static text values (so
# the saved sheet stores a literal string, not a formula) without
# disturbing the existing cell formatting/styles.
$dateRange = $ws.Range("A193:A195")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)
